$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15000
$ws.Range("I21").Value = 14444.444
$ws.Range("K21").Value = 14444.444
$ws.Range("M21").Value = -13976.444
$ws.Range("H23").Value = 15000
$ws.Range("I23").Value = 14444.444
$ws.Range("K23").Value = 14444.444
$ws.Range("M23").Value = -14210.444
$ws.Range("H62").Value = 31181.63
$ws.Range("I62").Value = 3704.6667
$ws.Range("J62").Value = 78285
$ws.Range("K62").Value = 3704.6667
$ws.Range("L62").Value = 78285
$ws.Range("M62").Value = -3080.6667
$ws.Range("N62").Value = -79533
$ws.Range("H65").Value = 31181.63
$ws.Range("I65").Value = 3704.6667
$ws.Range("J65").Value = 78285
$ws.Range("K65").Value = 18523.3335
$ws.Range("L65").Value = 391425
$ws.Range("M65").Value = -15403.3335
$ws.Range("N65").Value = -397665
$ws.Range("H70").Value = 602145.6
$ws.Range("I70").Value = 929614.6
$ws.Range("J70").Value = 1785.8334
$ws.Range("K70").Value = 2788843.8
$ws.Range("L70").Value = 5357.5002
$ws.Range("M70").Value = -2788573.8
$ws.Range("N70").Value = -5897.5002
$ws.Range("H73").Value = 602145.6
$ws.Range("I73").Value = 929614.6
$ws.Range("J73").Value = 1785.8334
$ws.Range("K73").Value = 2788843.8
$ws.Range("L73").Value = 5357.5002
$ws.Range("M73").Value = -2787907.8
$ws.Range("N73").Value = -7229.5002
$ws.Range("H98").Value = 1762.7084
$ws.Range("I98").Value = 1752.3914
$ws.Range("K98").Value = 1752.3914
$ws.Range("M98").Value = -254.3914
$ws.Range("H122").Value = 1762.7084
$ws.Range("I122").Value = 1752.3914
$ws.Range("K122").Value = 5257.174199999999
$ws.Range("M122").Value = -2807.174199999999
$ws.Range("H137").Value = 3080.5217
$ws.Range("I137").Value = 2397.625
$ws.Range("K137").Value = 7192.875
$ws.Range("M137").Value = -4642.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3486.9473
$ws.Range("I32").Value = 2402.9167
$ws.Range("K32").Value = 2402.9167
$ws.Range("M32").Value = -2115.9167
$ws.Range("H97").Value = 755.775
$ws.Range("I97").Value = 703.75
$ws.Range("J97").Value = 963.875
$ws.Range("K97").Value = 703.75
$ws.Range("L97").Value = 963.875
$ws.Range("M97").Value = -207.75
$ws.Range("N97").Value = -1955.875
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H132").Value = 2349.3953
$ws.Range("I132").Value = 1754.8379
$ws.Range("K132").Value = 5264.5137
$ws.Range("M132").Value = -2734.5137

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 16958.857
$ws.Range("I5").Value = 708.625
$ws.Range("J5").Value = 38625.832
$ws.Range("K5").Value = 708.625
$ws.Range("L5").Value = 38625.832
$ws.Range("M5").Value = -595.625
$ws.Range("N5").Value = -38851.832
$ws.Range("H7").Value = 200002300
$ws.Range("I7").Value = 252
$ws.Range("J7").Value = 333337000
$ws.Range("K7").Value = 252
$ws.Range("L7").Value = 333337000
$ws.Range("M7").Value = -139
$ws.Range("N7").Value = -333337226
$ws.Range("H20").Value = 9339.406000000001
$ws.Range("I20").Value = 10457.808
$ws.Range("K20").Value = 10457.808
$ws.Range("M20").Value = -10210.808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4347.533
$ws.Range("I132").Value = 4291.5415
$ws.Range("J132").Value = 4571.5
$ws.Range("K132").Value = 12874.6245
$ws.Range("L132").Value = 13714.5
$ws.Range("M132").Value = -10344.6245
$ws.Range("N132").Value = -18774.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 122588.15
$ws.Range("I11").Value = 6041.75
$ws.Range("J11").Value = 402299.5
$ws.Range("K11").Value = 18125.25
$ws.Range("L11").Value = 1206898.5
$ws.Range("M11").Value = -17985.25
$ws.Range("N11").Value = -1207178.5
$ws.Range("H64").Value = 550627.75
$ws.Range("J64").Value = 334166.34
$ws.Range("L64").Value = 1002499.02
$ws.Range("N64").Value = -1003039.02
$ws.Range("H67").Value = 550627.75
$ws.Range("J67").Value = 334166.34
$ws.Range("L67").Value = 1002499.02
$ws.Range("N67").Value = -1004371.02
$ws.Range("H70").Value = 8799.799999999999
$ws.Range("I70").Value = 6999.5
$ws.Range("K70").Value = 20998.5
$ws.Range("M70").Value = -20683.5
$ws.Range("H73").Value = 8799.799999999999
$ws.Range("I73").Value = 6999.5
$ws.Range("K73").Value = 20998.5
$ws.Range("M73").Value = -19906.5
$ws.Range("H80").Value = 4616.6665
$ws.Range("J80").Value = 4925
$ws.Range("L80").Value = 14775
$ws.Range("N80").Value = -16647
$ws.Range("H83").Value = 4616.6665
$ws.Range("J83").Value = 4925
$ws.Range("L83").Value = 44325
$ws.Range("N83").Value = -53685
$ws.Range("H92").Value = 480.16666
$ws.Range("I92").Value = 429.66666
$ws.Range("K92").Value = 1288.99998
$ws.Range("M92").Value = -40.99998000000005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4966.423
$ws.Range("I132").Value = 3926.5264
$ws.Range("J132").Value = 7789
$ws.Range("K132").Value = 11779.5792
$ws.Range("L132").Value = 23367
$ws.Range("M132").Value = -9249.5792
$ws.Range("N132").Value = -28427

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26318688
$ws.Range("I7").Value = 38464270
$ws.Range("K7").Value = 38464270
$ws.Range("M7").Value = -38464158
$ws.Range("H22").Value = 1481.3334
$ws.Range("I22").Value = 1496.5
$ws.Range("J22").Value = 1473.75
$ws.Range("K22").Value = 1496.5
$ws.Range("L22").Value = 1473.75
$ws.Range("M22").Value = -1201.5
$ws.Range("N22").Value = -2063.75
$ws.Range("H27").Value = 1481.3334
$ws.Range("I27").Value = 1496.5
$ws.Range("J27").Value = 1473.75
$ws.Range("K27").Value = 1496.5
$ws.Range("L27").Value = 1473.75
$ws.Range("M27").Value = -1389.5
$ws.Range("N27").Value = -1687.75
$ws.Range("H35").Value = 100001610
$ws.Range("I35").Value = 1790.2222
$ws.Range("K35").Value = 1790.2222
$ws.Range("M35").Value = -1454.2222
$ws.Range("H40").Value = 2639.5334
$ws.Range("I40").Value = 2650.2693
$ws.Range("J40").Value = 2569.75
$ws.Range("K40").Value = 2650.2693
$ws.Range("L40").Value = 2569.75
$ws.Range("M40").Value = -2514.2693
$ws.Range("N40").Value = -2841.75
$ws.Range("H51").Value = 43000
$ws.Range("J51").Value = 43000
$ws.Range("L51").Value = 43000
$ws.Range("N51").Value = -43956
$ws.Range("H55").Value = 681.3
$ws.Range("I55").Value = 818.75
$ws.Range("J55").Value = 589.6667
$ws.Range("K55").Value = 818.75
$ws.Range("L55").Value = 589.6667
$ws.Range("M55").Value = -645.75
$ws.Range("N55").Value = -935.6667
$ws.Range("H68").Value = 8171.727
$ws.Range("I68").Value = 4704.154
$ws.Range("J68").Value = 13180.444
$ws.Range("K68").Value = 4704.154
$ws.Range("L68").Value = 13180.444
$ws.Range("M68").Value = -3955.154
$ws.Range("N68").Value = -14678.444
$ws.Range("H69").Value = 50000
$ws.Range("J69").Value = 50000
$ws.Range("L69").Value = 50000
$ws.Range("N69").Value = -51622
$ws.Range("H71").Value = 8171.727
$ws.Range("I71").Value = 4704.154
$ws.Range("J71").Value = 13180.444
$ws.Range("K71").Value = 23520.77
$ws.Range("L71").Value = 65902.22
$ws.Range("M71").Value = -19776.77
$ws.Range("N71").Value = -73390.22
$ws.Range("H72").Value = 50000
$ws.Range("J72").Value = 50000
$ws.Range("L72").Value = 150000
$ws.Range("N72").Value = -158112
$ws.Range("H122").Value = 6052.2104
$ws.Range("I122").Value = 3874.3333
$ws.Range("K122").Value = 11622.9999
$ws.Range("M122").Value = -9172.999899999999
$ws.Range("H126").Value = 26318688
$ws.Range("I126").Value = 38464270
$ws.Range("K126").Value = 115392810
$ws.Range("M126").Value = -115390340
$ws.Range("H132").Value = 20003584
$ws.Range("I132").Value = 28573374
$ws.Range("J132").Value = 7407.2
$ws.Range("K132").Value = 85720122
$ws.Range("L132").Value = 22221.6
$ws.Range("M132").Value = -85717592
$ws.Range("N132").Value = -27281.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2009.9474
$ws.Range("I122").Value = 1399.5834
$ws.Range("J122").Value = 3056.2856
$ws.Range("K122").Value = 4198.7502
$ws.Range("L122").Value = 9168.856800000001
$ws.Range("M122").Value = -1748.7502
$ws.Range("N122").Value = -14068.8568
$ws.Range("H126").Value = 2021.85
$ws.Range("I126").Value = 1627.9231
$ws.Range("J126").Value = 2753.4285
$ws.Range("K126").Value = 4883.7693
$ws.Range("L126").Value = 8260.2855
$ws.Range("M126").Value = -2413.7693
$ws.Range("N126").Value = -13200.2855
$ws.Range("H132").Value = 4634.095
$ws.Range("I132").Value = 3436.742
$ws.Range("J132").Value = 8008.4546
$ws.Range("K132").Value = 10310.226
$ws.Range("L132").Value = 24025.3638
$ws.Range("M132").Value = -7780.226000000001
$ws.Range("N132").Value = -29085.3638

Write-Host "Applied all Diabolos_Profits market-data updates."
